$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "mapsto" column (column C); values shift left is not needed here,
# we just drop the whole column entirely so only mapsfrom/mapsto(meaning->renamed) remain.
$ws.Range("C1:C7").EntireColumn.Delete()

# Rename header B1 from "meaning" to "mapsto"
$ws.Range("B1").Value = "mapsto"

# Update / lowercase the stage names in column B
$ws.Range("B2").Value = "wake"
$ws.Range("B3").Value = "rem"
$ws.Range("B4").Value = "stage1"
$ws.Range("B5").Value = "stage2"
$ws.Range("B6").Value = "sws"
$ws.Range("B7").Value = "unknown"

# Update the active selection as recorded in the saved view state
$ws.Range("C5").Select()
